$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 37 (shifts old rows 37..121 down to 38..122)
$ws.Range("A37").EntireRow.Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A37").Value = 7
$ws.Range("B37").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C37").Value = 'Ñuble'
$ws.Range("D37").Value = 45177
$ws.Range("E37").Value = 16
$ws.Range("F37").Value = 100112001
$ws.Range("G37").Value = 'Berenjena'
$ws.Range("H37").Value = 'Sin especificar'
$ws.Range("I37").Value = 'Primera'
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 9000
$ws.Range("L37").Value = 9000
$ws.Range("M37").Value = 9000
$ws.Range("N37").Value = '$/caja 60 unidades'
$ws.Range("O37").Value = 'Región de Arica y Parinacota'
$ws.Range("P37").Value = 150
$ws.Range("Q37").Value = 60
$ws.Range("R37").Value = 'Hortaliza'
